$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 649.6
$ws.Range("I12").Value = 649.5
$ws.Range("J12").Value = 649.6667
$ws.Range("K12").Value = 649.5
$ws.Range("L12").Value = 649.6667
$ws.Range("M12").Value = -479.5
$ws.Range("N12").Value = -989.6667

$ws.Range("H33").Value = 120.25
$ws.Range("I33").Value = 99.125
$ws.Range("J33").Value = 162.5
$ws.Range("K33").Value = 99.125
$ws.Range("L33").Value = 162.5
$ws.Range("M33").Value = 129.875

$ws.Range("H38").Value = 62
$ws.Range("I38").Value = 63.5
$ws.Range("J38").Value = 59
$ws.Range("K38").Value = 190.5
$ws.Range("L38").Value = 177
$ws.Range("M38").Value = 181.5
$ws.Range("N38").Value = -921

$ws.Range("H39").Value = 237.78947
$ws.Range("I39").Value = 80.07692
$ws.Range("J39").Value = 579.5
$ws.Range("K39").Value = 240.23076
$ws.Range("L39").Value = 1738.5
$ws.Range("M39").Value = 55.76924
$ws.Range("N39").Value = -2330.5

$ws.Range("H62").Value = 2155
$ws.Range("I62").Value = 2155
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2155
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1531

$ws.Range("H65").Value = 2155
$ws.Range("I65").Value = 2155
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 10775
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -7655

$ws.Range("H127").Value = 1074.1428
$ws.Range("I127").Value = 1074.1428
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 3222.4284
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 1737.5716

$ws.Range("H138").Value = 2391.9333
$ws.Range("I138").Value = 871.8333
$ws.Range("J138").Value = 8472.333000000001
$ws.Range("K138").Value = 2615.4999
$ws.Range("L138").Value = 25416.999
$ws.Range("M138").Value = 2524.5001
$ws.Range("N138").Value = -35696.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1738.8
$ws.Range("I110").Value = 1576.8889
$ws.Range("J110").Value = 1981.6666
$ws.Range("K110").Value = 1576.8889
$ws.Range("L110").Value = 1981.6666
$ws.Range("M110").Value = 468.1111000000001

$ws.Range("H122").Value = 1393.3572
$ws.Range("I122").Value = 1393.3572
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4180.071599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1730.071599999999

$ws.Range("H128").Value = 119000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 119000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 119000
$ws.Range("N128").Value = -128960

$ws.Range("H131").Value = 45000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 45000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 45000
$ws.Range("N131").Value = -55080

$ws.Range("H132").Value = 2039.9375
$ws.Range("I132").Value = 2009.3334
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 6028.0002
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -3498.0002
$ws.Range("N132").Value = -12557

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1895.1765
$ws.Range("I86").Value = 1881.2667
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 1881.2667
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = -758.2666999999999

$ws.Range("H88").Value = 5932.5713
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5932.5713
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 5932.5713
$ws.Range("N88").Value = -6744.5713

$ws.Range("H89").Value = 1895.1765
$ws.Range("I89").Value = 1881.2667
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 9406.333499999999
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = -3790.333499999999

$ws.Range("H91").Value = 5932.5713
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5932.5713
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 5932.5713
$ws.Range("N91").Value = -8740.5713

$ws.Range("H99").Value = 1840
$ws.Range("I99").Value = 1736.6666
$ws.Range("J99").Value = 2150
$ws.Range("K99").Value = 1736.6666
$ws.Range("L99").Value = 2150
$ws.Range("M99").Value = -238.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1145.7142
$ws.Range("I16").Value = 1144.2
$ws.Range("J16").Value = 1149.5
$ws.Range("K16").Value = 1144.2
$ws.Range("L16").Value = 1149.5
$ws.Range("M16").Value = -857.2

$ws.Range("H31").Value = 4465.373
$ws.Range("I31").Value = 2525.9355
$ws.Range("J31").Value = 6135.4443
$ws.Range("K31").Value = 2525.9355
$ws.Range("L31").Value = 6135.4443
$ws.Range("M31").Value = -2230.9355

$ws.Range("H34").Value = 4465.373
$ws.Range("I34").Value = 2525.9355
$ws.Range("J34").Value = 6135.4443
$ws.Range("K34").Value = 2525.9355
$ws.Range("L34").Value = 6135.4443
$ws.Range("M34").Value = -2323.9355

$ws.Range("H113").Value = 1145.7142
$ws.Range("I113").Value = 1144.2
$ws.Range("J113").Value = 1149.5
$ws.Range("K113").Value = 1144.2
$ws.Range("L113").Value = 1149.5
$ws.Range("M113").Value = 1025.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2313911.2
$ws.Range("I4").Value = 2313911.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6941733.600000001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -6941621.600000001

$ws.Range("H55").Value = 7666.6665
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 7666.6665
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 22999.9995
$ws.Range("N55").Value = -23353.9995

$ws.Range("H68").Value = 700.6667
$ws.Range("I68").Value = 699
$ws.Range("J68").Value = 701.5
$ws.Range("K68").Value = 2097
$ws.Range("L68").Value = 2104.5
$ws.Range("M68").Value = -1286

$ws.Range("H71").Value = 700.6667
$ws.Range("I71").Value = 699
$ws.Range("J71").Value = 701.5
$ws.Range("K71").Value = 6291
$ws.Range("L71").Value = 6313.5
$ws.Range("M71").Value = -2235

$ws.Range("H103").Value = 1344.3572
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1344.3572
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 4033.0716
$ws.Range("N103").Value = -5791.071599999999

$ws.Range("H116").Value = 2309.6667
$ws.Range("I116").Value = 2309.6667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 6929.000100000001
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -3487.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 12089156
$ws.Range("I11").Value = 11767623
$ws.Range("J11").Value = 13000167
$ws.Range("K11").Value = 11767623
$ws.Range("L11").Value = 13000167
$ws.Range("M11").Value = -11767484
$ws.Range("N11").Value = -13000445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2841.5
$ws.Range("I40").Value = 2337.5
$ws.Range("J40").Value = 3849.5
$ws.Range("K40").Value = 2337.5
$ws.Range("L40").Value = 3849.5
$ws.Range("M40").Value = -2201.5

$ws.Range("H122").Value = 2809.6667
$ws.Range("I122").Value = 2604.125
$ws.Range("J122").Value = 2974.1
$ws.Range("K122").Value = 7812.375
$ws.Range("L122").Value = 8922.299999999999
$ws.Range("M122").Value = -5362.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1530
$ws.Range("I6").Value = 990
$ws.Range("J6").Value = 1800
$ws.Range("K6").Value = 990
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -875
$ws.Range("N6").Value = -2030

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -6376
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -31880
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 1995.2413
$ws.Range("I122").Value = 1588.0454
$ws.Range("J122").Value = 3275
$ws.Range("K122").Value = 4764.1362
$ws.Range("L122").Value = 9825
$ws.Range("M122").Value = -2314.1362
